# Mise à jour de l'application
# Fill in missing "Poids (kg)" (C) and "MG (%)" (D) measurements for the
# 1963-dated batch of rows (160-186), matching the style already used by
# the existing "MG (%)" cells (style of D2, numFmt 164 / "Pourcentage").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing D-column percentage format (style index already used
# throughout column D) so newly written D cells match exactly.
$ws.Range("D2").Copy()
$ws.Range("D160:D186").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New weight ("Poids") values in column C for rows that didn't have one yet.
$ws.Range("C162").Value = 76.7
$ws.Range("C167").Value = 64.3
$ws.Range("C168").Value = 70.8
$ws.Range("C174").Value = 60.9
$ws.Range("C175").Value = 76.5
$ws.Range("C176").Value = 72.6
$ws.Range("C177").Value = 78.5

# New body-fat percentage ("MG %") values in column D.
$ws.Range("D160").Value = 0.075
$ws.Range("D161").Value = 0.12
$ws.Range("D162").Value = 0.051
$ws.Range("D163").Value = 0.041
$ws.Range("D164").Value = 0.075
$ws.Range("D165").Value = 0.093
$ws.Range("D166").Value = 0.078
$ws.Range("D167").Value = 0.041
$ws.Range("D168").Value = 0.062
$ws.Range("D169").Value = 0.086
$ws.Range("D171").Value = 0.065
$ws.Range("D172").Value = 0.062
$ws.Range("D173").Value = 0.058
$ws.Range("D174").Value = 0.048
$ws.Range("D175").Value = 0.058
$ws.Range("D176").Value = 0.055
$ws.Range("D177").Value = 0.075
$ws.Range("D179").Value = 0.083
$ws.Range("D181").Value = 0.096
$ws.Range("D183").Value = 0.088
$ws.Range("D184").Value = 0.075
$ws.Range("D185").Value = 0.058
$ws.Range("D186").Value = 0.037
# D170, D178, D180, D182 stay blank (only the percentage style was applied).

# Restore the view state (scroll position / active cell) as left by the author.
$null = $ws.Range("I164").Select()
